# Rename header columns to concise key names (B, C, D, E, F, K, L, M, N, O)
# and drop the now-unneeded reviewer/tracking columns (P-T), after first
# preserving the "github repo" links (previously column Q) by moving them
# into column O ("link_to_source_code_optional" -> "github_link").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "tool_description"
$ws.Range("C1").Value = "contact_name"
$ws.Range("D1").Value = "contact_email"
$ws.Range("E1").Value = "is_maintained"
$ws.Range("F1").Value = "relevant_diseases"
$ws.Range("K1").Value = "required_expertise"
$ws.Range("L1").Value = "tool_type"
$ws.Range("M1").Value = "input_type"
$ws.Range("N1").Value = "docs_link"
$ws.Range("O1").Value = "github_link"

# Rows that carry a GitHub repo URL in column Q need that value moved into
# column O before the old reviewer/Q:T columns are removed.
$rowsWithGithubRepo = @(4, 6, 7, 10, 11, 13, 14, 21, 22, 24, 37, 40, 41, 43, 44)

foreach ($r in $rowsWithGithubRepo) {
    $qAddr = "Q" + $r
    $oAddr = "O" + $r
    $repoLink = $ws.Range($qAddr).Value2
    $ws.Range($oAddr).Value = $repoLink
}

# Remove the reviewer / github_repo / complete / pkg_dev_assessment /
# overall_assessment columns (P through T) now that the link has been kept.
$ws.Columns("P:T").Delete()
